$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 2 first (mirrors the author selecting the whole row before inserting
# a column), then insert a new column before column F. Everything from F..O
# shifts right to G..P, carrying over values/styles/validations/col widths.
$ws.Rows("2:2").Select() | Out-Null
$ws.Columns("F:F").Insert() | Out-Null

# New header for the inserted column, and matching input cell style below it
# (copy the "hyperlink-ish" input style used by the neighboring image columns).
$ws.Range("F1").Value = "VideoYouTube"
$ws.Range("I2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths: the new column F, and the merged G:H width used for the
# "Informacion adicional" / "Imagen Logo" columns after the shift.
$ws.Range("F1").ColumnWidth = 14.833333333333334
$ws.Range("G1:H1").ColumnWidth = 20.833333333333332

# Add a (non-restrictive) data validation on the whole new column, matching
# the "Diferenciador" message box used elsewhere in the sheet.
$dv = $ws.Range("F1:F1048576").Validation
$dv.Add(0, 1, 1)
$dv.ShowInput = $false
$dv.ShowError = $false
$dv.ErrorTitle = "Valor no válido"
$dv.ErrorMessage = "Por favor selecciona uno de los valores permitidos"
$dv.InputTitle = "Diferenciador"
$dv.InputMessage = "Selecciona un diferenciador"

# Re-select the full row 2 so the saved sheet view matches the post-insert state.
$ws.Rows("2:2").Select() | Out-Null
